$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Schedule")
$ws2 = $wb.Worksheets.Item("Detailed")

# Schedule sheet updates
$ws1.Range("E2").Value = 688.3251765
$ws1.Range("F2").Value = 11.38103797123016

# Detailed sheet updates (column B = Price, column C = Type)
$ws2.Range("B20").Value = 0.61183
$ws2.Range("B21").Value = -4.885
$ws2.Range("B22").Value = -3.81796

$ws2.Range("B23").Value = -2.83936
$ws2.Range("C23").Value = "historical"

$ws2.Range("C24").Value = "historical"

$ws2.Range("B25").Value = 0.00002
$ws2.Range("C25").Value = "historical"

$ws2.Range("B26").Value = -0.92111
$ws2.Range("C26").Value = "historical"

$ws2.Range("B27").Value = -0.92816
$ws2.Range("B28").Value = -5.50985
$ws2.Range("B29").Value = -5.50985
$ws2.Range("B30").Value = -2.65382
$ws2.Range("B31").Value = -0.94853
$ws2.Range("B32").Value = 0.00002
$ws2.Range("B33").Value = 0.00002
$ws2.Range("B34").Value = 0.51
$ws2.Range("B35").Value = 0.00005
$ws2.Range("B36").Value = 4.72713
$ws2.Range("B37").Value = 3.79328
$ws2.Range("B38").Value = 12.47395
$ws2.Range("B39").Value = 36.08921
$ws2.Range("B40").Value = 48.22136
$ws2.Range("B41").Value = 62.33686
$ws2.Range("B42").Value = 62.33685
$ws2.Range("B44").Value = 65
$ws2.Range("B46").Value = 64.8901
$ws2.Range("B47").Value = 64.53091999999999
$ws2.Range("B48").Value = 61.81891
$ws2.Range("B49").Value = 62.81918
